$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.814.40'
$ws.Range('D3').Value = '3.737.90'
$ws.Range('E3').Value = '  -1.77%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.30'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.73'
$ws.Range('E6').Value = '  -5.43%  '
$ws.Range('D7').Value = '3.737.65'
$ws.Range('E7').Value = '  -1.82%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +1.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.37'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.03'
$ws.Range('E13').Value = '  -0.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000244'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').Value = '4.361.74'
$ws.Range('E15').Value = '  -1.87%  '
$ws.Range('D16').Value = '3.735.84'
$ws.Range('E16').Value = '  -1.79%  '
$ws.Range('D17').Value = '68.776.16'
$ws.Range('E17').Value = '  +1.60%  '
$ws.Range('E18').Value = '  +0.84%  '
$ws.Range('E19').Value = '  +0.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.34'
$ws.Range('E20').Value = '  +5.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '497.32'
$ws.Range('E21').Value = '  +1.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.09'
$ws.Range('E22').Value = '  +11.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.724'
$ws.Range('E23').Value = '  -1.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.90'
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.31'
$ws.Range('E25').Value = '  -2.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000140'
$ws.Range('E26').Value = '  -7.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.37'
$ws.Range('E27').Value = '  +0.48%  '
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('E30').Value = '  -0.37%  '
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.91'
$ws.Range('E32').Value = '  +3.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.67'
$ws.Range('E33').Value = '  -3.27%  '
$ws.Range('D34').Value = '3.882.33'
$ws.Range('E34').Value = '  -1.74%  '
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('D36').Value = '3.669.06'
$ws.Range('E36').Value = '  -2.01%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.01'
$ws.Range('E38').Value = '  +0.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.81'
$ws.Range('E39').Value = '  +0.40%  '
$ws.Range('E40').Value = '  -1.93%  '
$ws.Range('E41').Value = '  -1.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '433.82'
$ws.Range('E42').Value = '  -3.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '49.06'
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.98'
$ws.Range('E44').Value = '  -1.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.87'
$ws.Range('E45').Value = '  -0.33%  '
$ws.Range('E46').Value = '  +1.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.65'
$ws.Range('E48').Value = '  -1.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.86'
$ws.Range('E49').Value = '  +1.44%  '
$ws.Range('E50').Value = '  +1.16%  '
$ws.Range('D51').Value = '2.744.26'
$ws.Range('E51').Value = '  -3.36%  '
